$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F2").Value = 94
$ws.Range("C3").Value = 78
$ws.Range("B9").Value = 377
$ws.Range("E9").Value = 471
$ws.Range("F9").Value = 541
$ws.Range("G9").Value = 434
$ws.Range("H9").Value = 442
$ws.Range("B10").Value = 1339
$ws.Range("C10").Value = 1592
$ws.Range("D10").Value = 1800
$ws.Range("E10").Value = 2185
$ws.Range("F10").Value = 2110
$ws.Range("G10").Value = 896
$ws.Range("H10").Value = 603
$ws.Range("I10").Value = 853
$ws.Range("J10").Value = 733
$ws.Range("B11").Value = 1849
$ws.Range("C11").Value = 2228
$ws.Range("D11").Value = 2450
$ws.Range("E11").Value = 2886
$ws.Range("F11").Value = 2895
$ws.Range("G11").Value = 1563
$ws.Range("H11").Value = 1327
$ws.Range("I11").Value = 1695
$ws.Range("J11").Value = 1536

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("B8").Value = 50
$ws.Range("E8").Value = 93
$ws.Range("G8").Value = 53
$ws.Range("B9").Value = 66
$ws.Range("E9").Value = 159
$ws.Range("G9").Value = 99

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F7").Value = 62
$ws.Range("B8").Value = 220
$ws.Range("C8").Value = 331
$ws.Range("D8").Value = 512
$ws.Range("E8").Value = 651
$ws.Range("F8").Value = 547
$ws.Range("I8").Value = 191
$ws.Range("B9").Value = 268
$ws.Range("C9").Value = 385
$ws.Range("D9").Value = 584
$ws.Range("E9").Value = 735
$ws.Range("F9").Value = 630
$ws.Range("I9").Value = 315

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("E7").Value = 39
$ws.Range("F7").Value = 55
$ws.Range("H7").Value = 7
$ws.Range("E8").Value = 52
$ws.Range("F8").Value = 93
$ws.Range("H8").Value = 24

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G8").Value = 91
$ws.Range("D10").Value = 32
$ws.Range("E24").Value = 2
$ws.Range("C27").Value = 26
$ws.Range("F28").Value = 123
$ws.Range("C30").Value = 26
$ws.Range("B32").Value = 66
$ws.Range("E32").Value = 159
$ws.Range("G32").Value = 99
$ws.Range("J35").Value = 18
$ws.Range("J42").Value = 22
$ws.Range("J47").Value = 43
$ws.Range("C49").Value = 15
$ws.Range("H52").Value = 17
$ws.Range("B53").Value = 268
$ws.Range("C53").Value = 385
$ws.Range("D53").Value = 584
$ws.Range("E53").Value = 735
$ws.Range("F53").Value = 630
$ws.Range("I53").Value = 315
$ws.Range("C62").Value = 28
$ws.Range("B63").Value = 16
$ws.Range("E65").Value = 52
$ws.Range("F65").Value = 93
$ws.Range("H65").Value = 24
$ws.Range("I70").Value = 34
$ws.Range("D72").Value = 13
$ws.Range("B74").Value = 48
$ws.Range("E76").Value = 95
$ws.Range("G76").Value = 47
$ws.Range("G77").Value = 54
$ws.Range("J81").Value = 11
$ws.Range("E89").Value = 22
$ws.Range("D92").Value = 42
$ws.Range("F98").Value = 11
$ws.Range("B99").Value = 1849
$ws.Range("C99").Value = 2228
$ws.Range("D99").Value = 2450
$ws.Range("E99").Value = 2886
$ws.Range("F99").Value = 2895
$ws.Range("G99").Value = 1563
$ws.Range("H99").Value = 1327
$ws.Range("I99").Value = 1695
$ws.Range("J99").Value = 1536

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("E4").Value = 5
$ws.Range("E6").Value = 22

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 11

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("F2").Value = 9
$ws.Range("F9").Value = 123

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J7").Value = 27
$ws.Range("J8").Value = 43

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("G8").Value = 9
$ws.Range("E9").Value = 76
$ws.Range("E10").Value = 95
$ws.Range("G10").Value = 47

$ws = $wb.Worksheets.Item('River North')
$ws.Range("B6").Value = 45
$ws.Range("B7").Value = 48

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("D8").Value = 34
$ws.Range("D9").Value = 42

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("H6").Value = 8
$ws.Range("H8").Value = 17

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("C6").Value = 20
$ws.Range("C7").Value = 26

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("C6").Value = 25
$ws.Range("C7").Value = 28

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("G9").Value = 30
$ws.Range("G10").Value = 54

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 15

$ws = $wb.Worksheets.Item('New City')
$ws.Range("B4").Value = 6
$ws.Range("B6").Value = 16

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 13

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 11

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J5").Value = 11
$ws.Range("J6").Value = 22

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("D6").Value = 29
$ws.Range("D7").Value = 32

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("C3").Value = 2
$ws.Range("C8").Value = 26

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("G7").Value = 57
$ws.Range("G8").Value = 91

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2
